$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.726.13'
$ws.Range("E2").Value = '  -0.27%  '
$ws.Range("D3").Value = '3.279.94'
$ws.Range("E3").Value = '  -0.77%  '
$ws.Range("D4").Value = '''0.999'
$ws.Range("E4").Value = '  -0.18%  '
$ws.Range("D5").Value = '''579.96'
$ws.Range("E5").Value = '  +4.08%  '
$ws.Range("D6").Value = '''183.80'
$ws.Range("E6").Value = '  -0.90%  '
$ws.Range("D7").Value = '''1.00'
$ws.Range("E7").Value = '  -0.05%  '
$ws.Range("D8").Value = '3.272.48'
$ws.Range("E8").Value = '  -0.76%  '
$ws.Range("D9").Value = '''0.570'
$ws.Range("E9").Value = '  -2.29%  '
$ws.Range("D10").Value = '''0.176'
$ws.Range("E10").Value = '  -5.04%  '
$ws.Range("E11").Value = '  -1.98%  '
$ws.Range("D12").Value = '''46.29'
$ws.Range("E12").Value = '  -2.43%  '
$ws.Range("E13").Value = '  -2.20%  '
$ws.Range("D14").Value = '3.807.07'
$ws.Range("E14").Value = '  -0.69%  '
$ws.Range("B15").Value = 'Polkadot'
$ws.Range("C15").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D15").Value = '''8.42'
$ws.Range("E15").Value = '  -2.48%  '
$ws.Range("B16").Value = 'BitcoinCash'
$ws.Range("C16").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D16").Value = '''618.35'
$ws.Range("E16").Value = '  -1.72%  '
$ws.Range("D17").Value = '65.676.97'
$ws.Range("E17").Value = '  -0.34%  '
$ws.Range("E18").Value = '  +0.42%  '
$ws.Range("D19").Value = '''17.81'
$ws.Range("E19").Value = '  -1.59%  '
$ws.Range("D20").Value = '3.284.10'
$ws.Range("E20").Value = '  -0.43%  '
$ws.Range("D21").Value = '''10.94'
$ws.Range("E21").Value = '  -3.53%  '
$ws.Range("D22").Value = '''0.889'
$ws.Range("E22").Value = '  -1.87%  '
$ws.Range("D23").Value = '''17.99'
$ws.Range("E23").Value = '  -0.24%  '
$ws.Range("D24").Value = '''100.97'
$ws.Range("E24").Value = '  -1.41%  '
$ws.Range("D25").Value = '''4.96'
$ws.Range("E25").Value = '  +0.31%  '
$ws.Range("D26").Value = '''4.02'
$ws.Range("E26").Value = '  +1.87%  '
$ws.Range("D27").Value = '''2.71'
$ws.Range("E27").Value = '  -0.10%  '
$ws.Range("E28").Value = '  -0.08%  '
$ws.Range("D29").Value = '''30.89'
$ws.Range("E29").Value = '  +2.32%  '
$ws.Range("D30").Value = '''8.42'
$ws.Range("E30").Value = '  -2.53%  '
$ws.Range("D31").Value = '''6.43'
$ws.Range("E31").Value = '  +0.70%  '
$ws.Range("D32").Value = '''3.76'
$ws.Range("E32").Value = '  -7.35%  '
$ws.Range("D33").Value = '''547.36'
$ws.Range("E33").Value = '  -0.12%  '
$ws.Range("D34").Value = '''10.85'
$ws.Range("E34").Value = '  -2.25%  '
$ws.Range("D35").Value = '3.786.70'
$ws.Range("E35").Value = '  -0.23%  '
$ws.Range("D36").Value = '''0.104'
$ws.Range("E36").Value = '  -1.44%  '
$ws.Range("D37").Value = '''1.00'
$ws.Range("E37").Value = '  -0.04%  '
$ws.Range("D38").Value = '''56.02'
$ws.Range("E38").Value = '  -2.54%  '
$ws.Range("D39").Value = '''0.128'
$ws.Range("E39").Value = '  -0.54%  '
$ws.Range("D40").Value = '''32.56'
$ws.Range("E40").Value = '  -3.10%  '
$ws.Range("B41").Value = 'Stacks'
$ws.Range("C41").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D41").Value = '''3.14'
$ws.Range("E41").Value = '  -3.32%  '
$ws.Range("B42").Value = 'ApeXProtocol'
$ws.Range("C42").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D42").Value = '''3.38'
$ws.Range("E42").Value = '  +2.89%  '
$ws.Range("D43").Value = '''2.58'
$ws.Range("E43").Value = '  -4.03%  '
$ws.Range("D44").Value = '0.0₃0678'
$ws.Range("E44").Value = '  -8.41%  '
$ws.Range("E45").Value = '  -1.22%  '
$ws.Range("D46").Value = '''0.0406'
$ws.Range("E46").Value = '  -2.53%  '
$ws.Range("D47").Value = '''2.98'
$ws.Range("E47").Value = '  -5.92%  '
$ws.Range("E48").Value = '  +0.25%  '
$ws.Range("E49").Value = '  -1.99%  '
$ws.Range("E50").Value = '  -3.67%  '
$ws.Range("D51").Value = '''128.24'
$ws.Range("E51").Value = '  +4.97%  '
